$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (shifts existing rows 17-25 down to 18-26)
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with Netherlands data
$ws.Range("A17").Value = "Netherlands"
$ws.Range("B17").Value = 339650
$ws.Range("C17").Value = 28495
$ws.Range("D17").Value = 33017
$ws.Range("E17").Value = 45500
$ws.Range("F17").Value = 29213
$ws.Range("G17").Value = 17969
$ws.Range("H17").Value = 14175
$ws.Range("I17").Value = 17407606
$ws.Range("J17").Value = 17475445
$ws.Range("K17").Value = 17441526

# Apply the number-format style (same as the rest of the numeric columns) to B17:K17
$ws.Range("B17:K17").NumberFormat = "#,##0"

# Update selection to match final state
$ws.Range("H17").Select()
